# Auto-generated Excel COM-interop script to apply Hades_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 405.55554
$ws.Range("I4").Value = 405.55554
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 405.55554
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -291.55554

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 460.7097
$ws.Range("I33").Value = 153
$ws.Range("J33").Value = 1106.9
$ws.Range("K33").Value = 153
$ws.Range("L33").Value = 1106.9
$ws.Range("M33").Value = 76
$ws.Range("N33").Value = -1564.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 16887.889
$ws.Range("I135").Value = 20958.715
$ws.Range("J135").Value = 2640
$ws.Range("K135").Value = 188628.435
$ws.Range("L135").Value = 23760
$ws.Range("M135").Value = -186093.435
$ws.Range("N135").Value = -28830

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2719809.5
$ws.Range("I138").Value = 209811.2
$ws.Range("J138").Value = 9806863
$ws.Range("K138").Value = 629433.6000000001
$ws.Range("L138").Value = 29420589
$ws.Range("M138").Value = -624293.6000000001
$ws.Range("N138").Value = -29430869

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 676.85565
$ws.Range("I32").Value = 691.92554
$ws.Range("J32").Value = 204.66667
$ws.Range("K32").Value = 691.92554
$ws.Range("L32").Value = 204.66667
$ws.Range("M32").Value = -404.92554
$ws.Range("N32").Value = -778.6666700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1461.8462
$ws.Range("I45").Value = 1452.7222
$ws.Range("J45").Value = 1482.375
$ws.Range("K45").Value = 1452.7222
$ws.Range("L45").Value = 1482.375
$ws.Range("M45").Value = -1075.7222
$ws.Range("N45").Value = -2236.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3269826.8
$ws.Range("I122").Value = 1907.0714
$ws.Range("J122").Value = 18520118
$ws.Range("K122").Value = 5721.2142
$ws.Range("L122").Value = 55560354
$ws.Range("M122").Value = -3271.2142
$ws.Range("N122").Value = -55565254

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 28000
$ws.Range("J55").Value = 28000
$ws.Range("L55").Value = 28000
$ws.Range("N55").Value = -28546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1838.1364
$ws.Range("I107").Value = 2152
$ws.Range("J107").Value = 1461.5
$ws.Range("K107").Value = 2152
$ws.Range("L107").Value = 1461.5
$ws.Range("M107").Value = -232
$ws.Range("N107").Value = -5301.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4375.923
$ws.Range("I31").Value = 2524.9092
$ws.Range("J31").Value = 5733.3335
$ws.Range("K31").Value = 2524.9092
$ws.Range("L31").Value = 5733.3335
$ws.Range("M31").Value = -2229.9092
$ws.Range("N31").Value = -6323.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4375.923
$ws.Range("I34").Value = 2524.9092
$ws.Range("J34").Value = 5733.3335
$ws.Range("K34").Value = 2524.9092
$ws.Range("L34").Value = 5733.3335
$ws.Range("M34").Value = -2322.9092
$ws.Range("N34").Value = -6137.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2065.8064
$ws.Range("I86").Value = 1999.381
$ws.Range("J86").Value = 2205.3
$ws.Range("K86").Value = 1999.381
$ws.Range("L86").Value = 2205.3
$ws.Range("M86").Value = -876.3810000000001
$ws.Range("N86").Value = -4451.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2065.8064
$ws.Range("I89").Value = 1999.381
$ws.Range("J89").Value = 2205.3
$ws.Range("K89").Value = 9996.905000000001
$ws.Range("L89").Value = 11026.5
$ws.Range("M89").Value = -4380.905000000001
$ws.Range("N89").Value = -22258.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6040.7144
$ws.Range("I99").Value = 6099.4287
$ws.Range("J99").Value = 5982
$ws.Range("K99").Value = 6099.4287
$ws.Range("L99").Value = 5982
$ws.Range("M99").Value = -4601.4287
$ws.Range("N99").Value = -8978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6040.7144
$ws.Range("I126").Value = 6099.4287
$ws.Range("J126").Value = 5982
$ws.Range("K126").Value = 18298.2861
$ws.Range("L126").Value = 17946
$ws.Range("M126").Value = -15828.2861
$ws.Range("N126").Value = -22886

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 745.2
$ws.Range("J122").Value = 1451
$ws.Range("L122").Value = 13059
$ws.Range("N122").Value = -17959

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 831.4286
$ws.Range("I132").Value = 787.5
$ws.Range("J132").Value = 890
$ws.Range("K132").Value = 7087.5
$ws.Range("L132").Value = 8010
$ws.Range("M132").Value = -4557.5
$ws.Range("N132").Value = -13070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3310.3333
$ws.Range("I31").Value = 1572.4
$ws.Range("J31").Value = 12000
$ws.Range("K31").Value = 1572.4
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = -1280.4
$ws.Range("N31").Value = -12584

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 3310.3333
$ws.Range("I37").Value = 1572.4
$ws.Range("J37").Value = 12000
$ws.Range("K37").Value = 1572.4
$ws.Range("L37").Value = 12000
$ws.Range("M37").Value = -1295.4
$ws.Range("N37").Value = -12554

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 113782.445
$ws.Range("I132").Value = 93142.91
$ws.Range("J132").Value = 146216
$ws.Range("K132").Value = 279428.73
$ws.Range("L132").Value = 438648
$ws.Range("M132").Value = -276898.73
$ws.Range("N132").Value = -443708

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 24298.344
$ws.Range("I132").Value = 12937.911
$ws.Range("J132").Value = 47535.59
$ws.Range("K132").Value = 38813.733
$ws.Range("L132").Value = 142606.77
$ws.Range("M132").Value = -36283.733
$ws.Range("N132").Value = -147666.77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 84096.16
$ws.Range("I136").Value = 51486.363
$ws.Range("J136").Value = 263450
$ws.Range("K136").Value = 154459.089
$ws.Range("L136").Value = 790350
$ws.Range("M136").Value = -151909.089
$ws.Range("N136").Value = -795450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10725
$ws.Range("I29").Value = 3000
$ws.Range("J29").Value = 13300
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 13300
$ws.Range("M29").Value = -2710
$ws.Range("N29").Value = -13880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 84571.46000000001
$ws.Range("I132").Value = 71883.61
$ws.Range("J132").Value = 124044.78
$ws.Range("K132").Value = 215650.83
$ws.Range("L132").Value = 372134.34
$ws.Range("M132").Value = -213120.83
$ws.Range("N132").Value = -377194.34

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 32944.906
$ws.Range("I136").Value = 23934.861
$ws.Range("K136").Value = 71804.583
$ws.Range("M136").Value = -69254.583
